$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aValues = @(45758.01041666666,45758.02083333334,45758.03125,45758.04166666666,45758.05208333334,45758.0625,45758.07291666666,45758.08333333334,45758.09375,45758.10416666666,45758.11458333334,45758.125,45758.13541666666,45758.14583333334,45758.15625,45758.16666666666,45758.17708333334,45758.1875,45758.19791666666,45758.20833333334,45758.21875,45758.22916666666,45758.23958333334,45758.25,45758.26041666666,45758.27083333334,45758.28125,45758.29166666666,45758.30208333334,45758.3125,45758.32291666666,45758.33333333334,45758.34375,45758.35416666666,45758.36458333334,45758.375,45758.38541666666,45758.39583333334,45758.40625,45758.41666666666,45758.42708333334,45758.4375,45758.44791666666,45758.45833333334,45758.46875,45758.47916666666,45758.48958333334,45758.5,45758.51041666666,45758.52083333334,45758.53125,45758.54166666666,45758.55208333334,45758.5625,45758.57291666666,45758.58333333334,45758.59375,45758.60416666666,45758.61458333334,45758.625,45758.63541666666,45758.64583333334,45758.65625,45758.66666666666,45758.67708333334,45758.6875,45758.69791666666,45758.70833333334,45758.71875,45758.72916666666,45758.73958333334,45758.75,45758.76041666666,45758.77083333334,45758.78125,45758.79166666666,45758.80208333334,45758.8125,45758.82291666666,45758.83333333334,45758.84375,45758.85416666666,45758.86458333334,45758.875,45758.88541666666,45758.89583333334,45758.90625,45758.91666666666,45758.92708333334,45758.9375,45758.94791666666,45758.95833333334,45758.96875,45758.97916666666,45758.98958333334,45759)
$bValues = @(1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,3,3,4,6,64,73,85,98,343,368,387,419,743,759,799,827,1040,1062,1081,1087,1200,1211,1218,1222,1258,1255,1253,1251,1169,1166,1160,1154,1080,1068,1053,1035,893,872,849,827,684,659,638,618,436,415,397,378,191,177,163,152,22,17,14,13,2,2,2,2,1,1,1,1,1,1,1,1,1,1,1,1,0,0,0,0)

$startRow = 2
for ($i = 0; $i -lt $aValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $aValues[$i]
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}
